$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain numeric-looking text (e.g. thousand-grouped "57.001.44").
# Force Text format on each changed Price cell first so Excel keeps the literal string
# instead of re-parsing it as a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '56.971.48'
$ws.Range("E2").Value = '  +6.97%  '

$ws.Range("D3").Value = '3.239.17'
$ws.Range("E3").Value = '  +2.70%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '394.55'
$ws.Range("E5").Value = '  -0.72%  '

$ws.Range("D6").Value = '107.09'
$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D7").Value = '3.235.25'
$ws.Range("E7").Value = '  +2.57%  '

$ws.Range("D8").Value = '0.563'
$ws.Range("E8").Value = '  +3.35%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").Value = '0.614'
$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("D11").Value = '38.83'
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").Value = '0.0980'
$ws.Range("E12").Value = '  +12.47%  '

$ws.Range("E13").Value = '  +1.80%  '

$ws.Range("D14").Value = '3.761.59'
$ws.Range("E14").Value = '  +3.11%  '

$ws.Range("D15").Value = '8.13'
$ws.Range("E15").Value = '  +1.35%  '

$ws.Range("D16").Value = '18.92'
$ws.Range("E16").Value = '  -0.47%  '

$ws.Range("D17").Value = '3.240.74'
$ws.Range("E17").Value = '  +2.68%  '

$ws.Range("E18").Value = '  -2.22%  '

$ws.Range("D19").Value = '10.74'
$ws.Range("E19").Value = '  -0.41%  '

$ws.Range("D20").Value = '56.810.91'
$ws.Range("E20").Value = '  +6.82%  '

$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("D22").Value = '0.0000105'
$ws.Range("E22").Value = '  +7.78%  '

$ws.Range("D23").Value = '12.98'
$ws.Range("E23").Value = '  +0.63%  '

$ws.Range("D24").Value = '296.88'
$ws.Range("E24").Value = '  +9.54%  '

$ws.Range("D25").Value = '73.49'
$ws.Range("E25").Value = '  +3.20%  '

$ws.Range("E26").Value = '  -2.71%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '27.88'
$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").Value = '4.38'
$ws.Range("E28").Value = '  +3.52%  '

$ws.Range("D29").Value = '7.72'
$ws.Range("E29").Value = '  -4.55%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  -4.32%  '

$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.168'
$ws.Range("E31").Value = '  -1.80%  '

$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("E33").Value = '  -0.57%  '

$ws.Range("D34").Value = '10.94'
$ws.Range("E34").Value = '  -0.87%  '

$ws.Range("D35").Value = '37.03'
$ws.Range("E35").Value = '  -0.98%  '

$ws.Range("E36").Value = '  -2.15%  '

$ws.Range("E37").Value = '  +1.46%  '

$ws.Range("D38").Value = '51.47'
$ws.Range("E38").Value = '  +1.91%  '

$ws.Range("D39").Value = '3.51'
$ws.Range("E39").Value = '  -0.32%  '

$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("D41").Value = '3.07'
$ws.Range("E41").Value = '  +11.23%  '

$ws.Range("D42").Value = '134.17'
$ws.Range("E42").Value = '  +3.40%  '

$ws.Range("D43").Value = '1.89'
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("E44").Value = '  +2.21%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '3.93'
$ws.Range("E45").Value = '  -5.36%  '

$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = '16.92'
$ws.Range("E46").Value = '  -2.42%  '

$ws.Range("D47").Value = '0.281'
$ws.Range("E47").Value = '  -4.17%  '

$ws.Range("D48").Value = '21.84'
$ws.Range("E48").Value = '  -2.38%  '

$ws.Range("D49").Value = '2.141.73'
$ws.Range("E49").Value = '  +2.46%  '

$ws.Range("E50").Value = '  -0.48%  '

$ws.Range("D51").Value = '2.03'
$ws.Range("E51").Value = '  +23.65%  '
